$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value2 = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1749153331348/OwQjboHWXD-camera-video-segment-driver-1749153333848.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNAQYB3VM%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T210029Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCaFLubTBgStMxFLXziVtR9XmN5AluFUp%2F2NI9xNNJVBgIgb%2FgY0V88iDSMwgHrtatlaohATlHYXJGHcM1t%2F1fw5Tsq3QMIZhAEGgw3ODEyMDQ5NDIyNDQiDBpCx0Ri5bazixVVTyq6AzkuckrLjLIbWM5NNdlgOTMgtArWz%2FhfCX0bMVfJ%2Bq%2BOsjCDM%2BDHx6Np5b%2BL9EZPaFeosBL04m67MOg5ELE%2F490ZHfLJlmjkjR%2BTLwwXkgB7oIeyey8v9UcAPWhl%2FF6UjVIXjnn8xJ1RBlqhkEAnGW2a9BCfZNJ1Mqg77%2FoX4F56knFRGhCwrl3XYWaxj8aeRWFjTHBB9%2F%2BoZV5QmuWn4dZTvJevy%2FEcOQxfGx0YP5RbVX9fssc2W%2Bp2uHBDMKf%2BeimadNBqwLq0NWnUOFamHl%2BvEY6MjnZCYL75nStouUYmv7u8q%2B4zxPAAcaVpBjzPnW1QvI6qgIpLCtoKaOz%2FtOiBwoJxhCE16YRRpJVVA9%2FD3FuSTEXHf5omOPG0kch5Coesiv9z%2Bs%2FTxaNXGuqU6yBzXIHnYKvF3xpoyaaEPQrH6Oc6ICKD9Ww1eE66rgZ%2BorcxI1rjUHlQxJ%2FLZ2fk35jM2XBozu%2BM3oGANiZnTQ85z%2BYT3XOAtrdDPB1BECw3QKhlOdlRYkcQjxuSEk94kIC9KU7v7TSAl9y%2Ff5yQ29qFol%2F2pVIjlaWmI6JAyOR37FT9GS2M0i%2BO65Mw%2FaiNwgY6pQE%2BNdVQFd2qzIbHTVepJpszOnlWHb2NLMKbzLqQNeDhWu8K09drnwEyGvRfEs001%2BME0Ubp061Oh8%2FIRUwYYCBRjngUhQw17h6fol846y8adlUvzZJHn4PShwfeeIgGWwkUYuVMj599KpJZ6xwLHrByxHJZ4PBIWAuoaENr3wiluKfNjghIh4ETZmUJerWnL73gg8WwuQ23QNZh%2F7gU2y%2BJ6pEufmc%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2005%3A00%3A29%20GMT&X-Amz-Signature=12d83c79f714dd8a87da54ba286f2c6f89d32f34b47cf578b6ac89eaf704c055'
$ws.Range("L3").Value2 = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1749153331348/oV5nVwJUOa-camera-video-segment-1749153333848.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNAQYB3VM%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T210029Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCaFLubTBgStMxFLXziVtR9XmN5AluFUp%2F2NI9xNNJVBgIgb%2FgY0V88iDSMwgHrtatlaohATlHYXJGHcM1t%2F1fw5Tsq3QMIZhAEGgw3ODEyMDQ5NDIyNDQiDBpCx0Ri5bazixVVTyq6AzkuckrLjLIbWM5NNdlgOTMgtArWz%2FhfCX0bMVfJ%2Bq%2BOsjCDM%2BDHx6Np5b%2BL9EZPaFeosBL04m67MOg5ELE%2F490ZHfLJlmjkjR%2BTLwwXkgB7oIeyey8v9UcAPWhl%2FF6UjVIXjnn8xJ1RBlqhkEAnGW2a9BCfZNJ1Mqg77%2FoX4F56knFRGhCwrl3XYWaxj8aeRWFjTHBB9%2F%2BoZV5QmuWn4dZTvJevy%2FEcOQxfGx0YP5RbVX9fssc2W%2Bp2uHBDMKf%2BeimadNBqwLq0NWnUOFamHl%2BvEY6MjnZCYL75nStouUYmv7u8q%2B4zxPAAcaVpBjzPnW1QvI6qgIpLCtoKaOz%2FtOiBwoJxhCE16YRRpJVVA9%2FD3FuSTEXHf5omOPG0kch5Coesiv9z%2Bs%2FTxaNXGuqU6yBzXIHnYKvF3xpoyaaEPQrH6Oc6ICKD9Ww1eE66rgZ%2BorcxI1rjUHlQxJ%2FLZ2fk35jM2XBozu%2BM3oGANiZnTQ85z%2BYT3XOAtrdDPB1BECw3QKhlOdlRYkcQjxuSEk94kIC9KU7v7TSAl9y%2Ff5yQ29qFol%2F2pVIjlaWmI6JAyOR37FT9GS2M0i%2BO65Mw%2FaiNwgY6pQE%2BNdVQFd2qzIbHTVepJpszOnlWHb2NLMKbzLqQNeDhWu8K09drnwEyGvRfEs001%2BME0Ubp061Oh8%2FIRUwYYCBRjngUhQw17h6fol846y8adlUvzZJHn4PShwfeeIgGWwkUYuVMj599KpJZ6xwLHrByxHJZ4PBIWAuoaENr3wiluKfNjghIh4ETZmUJerWnL73gg8WwuQ23QNZh%2F7gU2y%2BJ6pEufmc%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2005%3A00%3A29%20GMT&X-Amz-Signature=1d907d00f1202b1abbdd5d58567a855d605a9b6b1979247b35484ceb213075d7'
$ws.Range("K4").Value2 = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109254/1749147387468/GKgOu1ZeR3-camera-video-segment-driver-1749147389968.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNAQYB3VM%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T210029Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCaFLubTBgStMxFLXziVtR9XmN5AluFUp%2F2NI9xNNJVBgIgb%2FgY0V88iDSMwgHrtatlaohATlHYXJGHcM1t%2F1fw5Tsq3QMIZhAEGgw3ODEyMDQ5NDIyNDQiDBpCx0Ri5bazixVVTyq6AzkuckrLjLIbWM5NNdlgOTMgtArWz%2FhfCX0bMVfJ%2Bq%2BOsjCDM%2BDHx6Np5b%2BL9EZPaFeosBL04m67MOg5ELE%2F490ZHfLJlmjkjR%2BTLwwXkgB7oIeyey8v9UcAPWhl%2FF6UjVIXjnn8xJ1RBlqhkEAnGW2a9BCfZNJ1Mqg77%2FoX4F56knFRGhCwrl3XYWaxj8aeRWFjTHBB9%2F%2BoZV5QmuWn4dZTvJevy%2FEcOQxfGx0YP5RbVX9fssc2W%2Bp2uHBDMKf%2BeimadNBqwLq0NWnUOFamHl%2BvEY6MjnZCYL75nStouUYmv7u8q%2B4zxPAAcaVpBjzPnW1QvI6qgIpLCtoKaOz%2FtOiBwoJxhCE16YRRpJVVA9%2FD3FuSTEXHf5omOPG0kch5Coesiv9z%2Bs%2FTxaNXGuqU6yBzXIHnYKvF3xpoyaaEPQrH6Oc6ICKD9Ww1eE66rgZ%2BorcxI1rjUHlQxJ%2FLZ2fk35jM2XBozu%2BM3oGANiZnTQ85z%2BYT3XOAtrdDPB1BECw3QKhlOdlRYkcQjxuSEk94kIC9KU7v7TSAl9y%2Ff5yQ29qFol%2F2pVIjlaWmI6JAyOR37FT9GS2M0i%2BO65Mw%2FaiNwgY6pQE%2BNdVQFd2qzIbHTVepJpszOnlWHb2NLMKbzLqQNeDhWu8K09drnwEyGvRfEs001%2BME0Ubp061Oh8%2FIRUwYYCBRjngUhQw17h6fol846y8adlUvzZJHn4PShwfeeIgGWwkUYuVMj599KpJZ6xwLHrByxHJZ4PBIWAuoaENr3wiluKfNjghIh4ETZmUJerWnL73gg8WwuQ23QNZh%2F7gU2y%2BJ6pEufmc%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2005%3A00%3A29%20GMT&X-Amz-Signature=19cad96118bd4659ea93b8b7b86dacecb28e79b22496f315c1a7183479fb427d'
$ws.Range("L4").Value2 = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109254/1749147387468/dQoeD1EUOr-camera-video-segment-1749147389968.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNAQYB3VM%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T210029Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCaFLubTBgStMxFLXziVtR9XmN5AluFUp%2F2NI9xNNJVBgIgb%2FgY0V88iDSMwgHrtatlaohATlHYXJGHcM1t%2F1fw5Tsq3QMIZhAEGgw3ODEyMDQ5NDIyNDQiDBpCx0Ri5bazixVVTyq6AzkuckrLjLIbWM5NNdlgOTMgtArWz%2FhfCX0bMVfJ%2Bq%2BOsjCDM%2BDHx6Np5b%2BL9EZPaFeosBL04m67MOg5ELE%2F490ZHfLJlmjkjR%2BTLwwXkgB7oIeyey8v9UcAPWhl%2FF6UjVIXjnn8xJ1RBlqhkEAnGW2a9BCfZNJ1Mqg77%2FoX4F56knFRGhCwrl3XYWaxj8aeRWFjTHBB9%2F%2BoZV5QmuWn4dZTvJevy%2FEcOQxfGx0YP5RbVX9fssc2W%2Bp2uHBDMKf%2BeimadNBqwLq0NWnUOFamHl%2BvEY6MjnZCYL75nStouUYmv7u8q%2B4zxPAAcaVpBjzPnW1QvI6qgIpLCtoKaOz%2FtOiBwoJxhCE16YRRpJVVA9%2FD3FuSTEXHf5omOPG0kch5Coesiv9z%2Bs%2FTxaNXGuqU6yBzXIHnYKvF3xpoyaaEPQrH6Oc6ICKD9Ww1eE66rgZ%2BorcxI1rjUHlQxJ%2FLZ2fk35jM2XBozu%2BM3oGANiZnTQ85z%2BYT3XOAtrdDPB1BECw3QKhlOdlRYkcQjxuSEk94kIC9KU7v7TSAl9y%2Ff5yQ29qFol%2F2pVIjlaWmI6JAyOR37FT9GS2M0i%2BO65Mw%2FaiNwgY6pQE%2BNdVQFd2qzIbHTVepJpszOnlWHb2NLMKbzLqQNeDhWu8K09drnwEyGvRfEs001%2BME0Ubp061Oh8%2FIRUwYYCBRjngUhQw17h6fol846y8adlUvzZJHn4PShwfeeIgGWwkUYuVMj599KpJZ6xwLHrByxHJZ4PBIWAuoaENr3wiluKfNjghIh4ETZmUJerWnL73gg8WwuQ23QNZh%2F7gU2y%2BJ6pEufmc%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2005%3A00%3A29%20GMT&X-Amz-Signature=80d4ac962d2faf43a60df65aa06ec530fe8c461383a93f090f2864c9dbb8d160'
